$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 14-23: convert project_id (C) and epochs (K) columns from text to numeric ---
$numericFixRows = @{
    14 = @{ C = 21691136843814; K = 10 }
    15 = @{ C = 103610396583974; K = 10 }
    16 = @{ C = 21691136843814; K = 10 }
    17 = @{ C = 21691136843814; K = 10 }
    18 = @{ C = 1938830; K = 10 }
    19 = @{ C = 1938830; K = 10 }
    20 = @{ C = 102004060440613; K = 10 }
    21 = @{ C = 6464689526794; K = 10 }
    22 = @{ C = 6464689526794; K = 10 }
    23 = @{ C = 6464689526794; K = 10 }
}

foreach ($r in $numericFixRows.Keys) {
    $vals = $numericFixRows[$r]
    $ws.Range("C$r").Value = $vals.C
    $ws.Range("K$r").Value = $vals.K
}

# Row 15: dev/test columns (H15, I15) lose their "N/A" text, becoming blank cells
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = ""

# --- New rows 24-33 appended to the report ---
$newRows = @(
    @("Large Images", "segmentation",      "21691136843814",  "7.11 minutes", "SegFormer-[14M]",  "{'height': 6000, 'width': 6000, 'paddingValue': 0}", "72%", "39%", "55%", "2025-06-10 11:10:25", "10"),
    @("Large Images", "object-detection",  "103610396583974", "3.13 minutes", "RepPoints-[37M]",  "{'height': 2000, 'width': 2000, 'paddingValue': 0}", "83%", "N/A", "N/A", "2025-06-10 11:16:03", "10"),
    @("Large Images", "segmentation",      "21691136843814",  "6.96 minutes", "SegFormer-[14M]",  "{'height': 6000, 'width': 6000, 'paddingValue': 0}", "73%", "43%", "56%", "2025-06-10 11:25:28", "10"),
    @("Large Images", "segmentation",      "21691136843814",  "7.96 minutes", "FastVit-[14M]",    "{'height': 6000, 'width': 6000, 'paddingValue': 0}", "73%", "48%", "58%", "2025-06-10 11:36:08", "10"),
    @("Large Images", "segmentation",      "1938830",         "4.56 minutes", "FastVit-[14M]",    "{'height': 2048, 'width': 2048, 'paddingValue': 0}", "8%",  "10%", "4%",  "2025-06-10 11:42:55", "10"),
    @("Large Images", "segmentation",      "1938830",         "4.14 minutes", "SegFormer-[14M]",  "{'height': 2048, 'width': 2048, 'paddingValue': 0}", "29%", "11%", "43%", "2025-06-10 11:49:32", "10"),
    @("Large Images", "object-detection",  "102004060440613", "9.72 minutes", "RtmDet-[9M]",      "{'height': 3040, 'width': 4056, 'paddingValue': 0}", "66%", "76%", "71%", "2025-06-10 12:01:34", "10"),
    @("Large Images", "object-detection",  "6464689526794",   "6.42 minutes", "RtmDet-[9M]",      "{'height': 6000, 'width': 6000, 'paddingValue': 0}", "69%", "70%", "67%", "2025-06-10 12:10:09", "10"),
    @("Large Images", "object-detection",  "6464689526794",   "7.30 minutes", "RepPoints-[20M]",  "{'height': 6000, 'width': 6000, 'paddingValue': 0}", "77%", "76%", "73%", "2025-06-10 12:20:31", "10"),
    @("Large Images", "object-detection",  "6464689526794",   "9.72 minutes", "RepPoints-[37M]",  "{'height': 6000, 'width': 6000, 'paddingValue': 0}", "69%", "71%", "68%", "2025-06-10 12:33:25", "10")
)

$startRow = 24
# Columns whose text values look numeric/percentage and must be forced to
# stay text (matches the source file, which stores these as inline strings
# rather than numbers for the newly appended rows):
#   C = project_id ("21691136843814"), G/H/I = train/dev/test ("72%"...),
#   K = epochs ("10")
$textCols = @(3, 7, 8, 9, 11)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($textCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$c - 1]
    }
}
